$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 184: bva rolling (1/28)
$ws.Cells.Item(184,1).Value = 78
$ws.Cells.Item(184,2).Value = 2022
$ws.Cells.Item(184,3).Value = 22
$ws.Cells.Item(184,4).Value = 1
$ws.Cells.Item(184,5).Value = 25
$ws.Cells.Item(184,6).Value = "bva"
$ws.Cells.Item(184,7).Value = "online"
$ws.Cells.Item(184,8).Value = "partially"
$ws.Cells.Item(184,9).Value = 0
$ws.Cells.Item(184,10).Value = 910
$ws.Cells.Item(184,11).Value = 0.5
$ws.Cells.Item(184,12).Value = "T_0.5"
$ws.Cells.Item(184,13).Value = 10
$ws.Cells.Item(184,14).Value = 3
$ws.Cells.Item(184,16).Value = 6
$ws.Cells.Item(184,17).Value = 3
$ws.Cells.Item(184,18).Value = 24
$ws.Cells.Item(184,19).Value = 16
$ws.Cells.Item(184,22).Value = 0.5
$ws.Cells.Item(184,23).Value = 2
$ws.Cells.Item(184,24).Value = 18
$ws.Cells.Item(184,25).Value = 12.5
$ws.Cells.Item(184,26).Value = 0.5
$ws.Cells.Item(184,30).Value = 4
$ws.Cells.Item(184,31).Value = "T_0.5"

# Row 185: ifop rolling (1/28)
$ws.Cells.Item(185,1).Value = 79
$ws.Cells.Item(185,2).Value = 2022
$ws.Cells.Item(185,3).Value = 22
$ws.Cells.Item(185,4).Value = 1
$ws.Cells.Item(185,5).Value = 27
$ws.Cells.Item(185,6).Value = "ifop"
$ws.Cells.Item(185,7).Value = "online"
$ws.Cells.Item(185,8).Value = "included"
$ws.Cells.Item(185,9).Value = 1
$ws.Cells.Item(185,10).Value = 1000
$ws.Cells.Item(185,11).Value = 0.5
$ws.Cells.Item(185,12).Value = "T_0.5"
$ws.Cells.Item(185,13).Value = 9.5
$ws.Cells.Item(185,14).Value = 3
$ws.Cells.Item(185,16).Value = 5.5
$ws.Cells.Item(185,17).Value = 3.5
$ws.Cells.Item(185,18).Value = 24
$ws.Cells.Item(185,19).Value = 16.5
$ws.Cells.Item(185,22).Value = 1
$ws.Cells.Item(185,23).Value = 1
$ws.Cells.Item(185,24).Value = 18
$ws.Cells.Item(185,25).Value = 13.5
$ws.Cells.Item(185,30).Value = 4
$ws.Cells.Item(185,31).Value = "T_0.5"

# Row 186: opinionway rolling (1/25)
$ws.Cells.Item(186,1).Value = 80
$ws.Cells.Item(186,2).Value = 2022
$ws.Cells.Item(186,3).Value = 22
$ws.Cells.Item(186,4).Value = 1
$ws.Cells.Item(186,5).Value = 24
$ws.Cells.Item(186,6).Value = "opinionway"
$ws.Cells.Item(186,7).Value = "online"
$ws.Cells.Item(186,8).Value = "included"
$ws.Cells.Item(186,9).Value = 1
$ws.Cells.Item(186,10).Value = 1000
$ws.Cells.Item(186,11).Value = 1
$ws.Cells.Item(186,12).Value = "T_0.5"
$ws.Cells.Item(186,13).Value = 9
$ws.Cells.Item(186,14).Value = 3
$ws.Cells.Item(186,16).Value = 5
$ws.Cells.Item(186,17).Value = 3
$ws.Cells.Item(186,18).Value = 25
$ws.Cells.Item(186,19).Value = 17
$ws.Cells.Item(186,22).Value = 1
$ws.Cells.Item(186,23).Value = 1
$ws.Cells.Item(186,24).Value = 16
$ws.Cells.Item(186,25).Value = 14
$ws.Cells.Item(186,30).Value = 5

# Update the view: scroll to show column S, adjust frozen pane & selection
$ws.Range("S1").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A173").Select()
$win.FreezePanes = $true
$ws.Range("AF185").Select()
